# --- Part 1: "day" sheet - convert D58:D63 from text to numeric bsecode ---
$wb = $excel.ActiveWorkbook
$wsDay = $wb.Worksheets.Item("day")

$bsecodes = @(500488, 532500, 532478, 533309, 500570, 532733)
for ($i = 0; $i -lt $bsecodes.Length; $i++) {
    $wsDay.Cells.Item(58 + $i, 4).Value = $bsecodes[$i]
}

# --- Part 2: "week" sheet - append new rows 72:82 ---
$wsWeek = $wb.Worksheets.Item("week")

$newRows = @(
    @(72, 1, "ALKEM", "Alkem Laboratories Limited", "539523", 1.38, 4991.35, 272433, "week", "28/06/2024 11:32:25"),
    @(73, 2, "LUPIN", "Lupin Limited", "500257", 2.4, 1621.35, 1413997, "week", "28/06/2024 11:32:25"),
    @(74, 3, "SUNPHARMA", "Sun Pharmaceuticals Industries Limited", "524715", 0.3, 1520.85, 2634599, "week", "28/06/2024 11:32:25"),
    @(75, 4, "ZYDUSLIFE", "Zydus Lifesciences Ltd", "532321", 1.6, 1073.95, 2188481, "week", "28/06/2024 11:32:25"),
    @(76, 5, "IRCTC", "Indian Railway Catering & Tourism Corporation Ltd", "542830", -0.2, 989.25, 3777772, "week", "28/06/2024 11:32:25"),
    @(77, 6, "VEDL", "Vedanta Limited", "500295", 2.41, 454, 15248617, "week", "28/06/2024 11:32:25"),
    @(78, 7, "LAURUSLABS", "Laurus Labs Limited", "540222", 0.25, 424.55, 793710, "week", "28/06/2024 11:32:25"),
    @(79, 8, "HINDCOPPER", "Hindustan Copper Limited", "513599", -0.3, 318.7, 6389187, "week", "28/06/2024 11:32:25"),
    @(80, 9, "BHEL", "Bharat Heavy Electricals Limited", "500103", 1.31, 300.85, 32298955, "week", "28/06/2024 11:32:25"),
    @(81, 10, "NMDC", "Nmdc Limited", "526371", 0.22, 246.05, 12936396, "week", "28/06/2024 11:32:25"),
    @(82, 11, "SAIL", "Steel Authority Of India Limited", "500113", 4.04, 148.65, 72589456, "week", "28/06/2024 11:32:25")
)

# Column D (bsecode) keeps storing text, same as the other not-yet-normalised
# rows already present in this sheet -- format the whole new block as Text
# up front so the literal numeric-looking strings aren't auto-converted.
$dRange = $wsWeek.Range("D72:D82")
$dRange.NumberFormat = "@"

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $wsWeek.Cells.Item($rowNum, 1).Value = $r[1]
    $wsWeek.Cells.Item($rowNum, 2).Value = $r[2]
    $wsWeek.Cells.Item($rowNum, 3).Value = $r[3]
    $wsWeek.Cells.Item($rowNum, 4).Value = $r[4]
    $wsWeek.Cells.Item($rowNum, 5).Value = $r[5]
    $wsWeek.Cells.Item($rowNum, 6).Value = $r[6]
    $wsWeek.Cells.Item($rowNum, 7).Value = $r[7]
    $wsWeek.Cells.Item($rowNum, 8).Value = $r[8]
    $wsWeek.Cells.Item($rowNum, 9).Value = $r[9]
}

# Drop the temporary Text number format again so the new cells end up with
# the same (default) style as every other data row in the sheet.
$dRange.Style = "Normal"
